{"js": "// Replace the 25 \"three-digit number divided by one-digit number\" answer\n// cells with their new values (1-to-1, same order as they appear in the\n// document body).\nconst replacements = [\n  [\"725\u00f75=145, 0\", \"663\u00f79=73, 6\"],\n  [\"510\u00f75=102, 0\", \"283\u00f74=70, 3\"],\n  [\"610\u00f74=152, 2\", \"933\u00f72=466, 1\"],\n  [\"498\u00f79=55, 3\", \"703\u00f77=100, 3\"],\n  [\"261\u00f72=130, 1\", \"964\u00f78=120, 4\"],\n  [\"198\u00f74=49, 2\", \"918\u00f77=131, 1\"],\n  [\"165\u00f72=82, 1\", \"707\u00f72=353, 1\"],\n  [\"979\u00f76=163, 1\", \"261\u00f79=29, 0\"],\n  [\"230\u00f76=38, 2\", \"722\u00f79=80, 2\"],\n  [\"703\u00f74=175, 3\", \"935\u00f75=187, 0\"],\n  [\"882\u00f76=147, 0\", \"969\u00f78=121, 1\"],\n  [\"997\u00f73=332, 1\", \"125\u00f73=41, 2\"],\n  [\"146\u00f74=36, 2\", \"556\u00f76=92, 4\"],\n  [\"137\u00f77=19, 4\", \"958\u00f77=136, 6\"],\n  [\"559\u00f78=69, 7\", \"773\u00f74=193, 1\"],\n  [\"643\u00f76=107, 1\", \"423\u00f78=52, 7\"],\n  [\"284\u00f72=142, 0\", \"295\u00f73=98, 1\"],\n  [\"986\u00f79=109, 5\", \"839\u00f78=104, 7\"],\n  [\"981\u00f74=245, 1\", \"113\u00f73=37, 2\"],\n  [\"679\u00f79=75, 4\", \"606\u00f74=151, 2\"],\n  [\"317\u00f75=63, 2\", \"892\u00f72=446, 0\"],\n  [\"381\u00f72=190, 1\", \"319\u00f78=39, 7\"],\n  [\"327\u00f77=46, 5\", \"307\u00f72=153, 1\"],\n  [\"525\u00f75=105, 0\", \"181\u00f77=25, 6\"],\n  [\"393\u00f79=43, 6\", \"184\u00f76=30, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit number divided by one-digit number\" answer\n# cells with their new values (1-to-1, same order as they appear in the\n# document body).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"725\u00f75=145, 0\", \"663\u00f79=73, 6\"),\n    @(\"510\u00f75=102, 0\", \"283\u00f74=70, 3\"),\n    @(\"610\u00f74=152, 2\", \"933\u00f72=466, 1\"),\n    @(\"498\u00f79=55, 3\", \"703\u00f77=100, 3\"),\n    @(\"261\u00f72=130, 1\", \"964\u00f78=120, 4\"),\n    @(\"198\u00f74=49, 2\", \"918\u00f77=131, 1\"),\n    @(\"165\u00f72=82, 1\", \"707\u00f72=353, 1\"),\n    @(\"979\u00f76=163, 1\", \"261\u00f79=29, 0\"),\n    @(\"230\u00f76=38, 2\", \"722\u00f79=80, 2\"),\n    @(\"703\u00f74=175, 3\", \"935\u00f75=187, 0\"),\n    @(\"882\u00f76=147, 0\", \"969\u00f78=121, 1\"),\n    @(\"997\u00f73=332, 1\", \"125\u00f73=41, 2\"),\n    @(\"146\u00f74=36, 2\", \"556\u00f76=92, 4\"),\n    @(\"137\u00f77=19, 4\", \"958\u00f77=136, 6\"),\n    @(\"559\u00f78=69, 7\", \"773\u00f74=193, 1\"),\n    @(\"643\u00f76=107, 1\", \"423\u00f78=52, 7\"),\n    @(\"284\u00f72=142, 0\", \"295\u00f73=98, 1\"),\n    @(\"986\u00f79=109, 5\", \"839\u00f78=104, 7\"),\n    @(\"981\u00f74=245, 1\", \"113\u00f73=37, 2\"),\n    @(\"679\u00f79=75, 4\", \"606\u00f74=151, 2\"),\n    @(\"317\u00f75=63, 2\", \"892\u00f72=446, 0\"),\n    @(\"381\u00f72=190, 1\", \"319\u00f78=39, 7\"),\n    @(\"327\u00f77=46, 5\", \"307\u00f72=153, 1\"),\n    @(\"525\u00f75=105, 0\", \"181\u00f77=25, 6\"),\n    @(\"393\u00f79=43, 6\", \"184\u00f76=30, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
